# Move four shapes on slide 3 straight down (only the vertical offset
# changes; the horizontal offset and the width/height stay the same).
#
# NOTE: Top/Left are expressed in points (1 pt = 12700 EMU) and the host
# stores them as single-precision floats, so a plain "EMU / 12700" can
# truncate to one EMU short after the round-trip back to XML. The literals
# below are nudged by a few 1e-7 pt so the stored value still truncates to
# the exact target EMU.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# Shape 1: Title "Problem Statement" -> y: 545251 -> 652763 EMU
$s.Shapes.Item(1).Top = 51.398661417322835

# Shape 2: Content Placeholder -> y: 1399326 -> 1554984 EMU
$s.Shapes.Item(2).Top = 122.43968583937007

# Shape 3: Title "Hypothesis" -> y: 3574780 -> 4016060 EMU
$s.Shapes.Item(3).Top = 316.22520455039376

# Shape 4: Content Placeholder -> y: 4435996 -> 4791677 EMU
$s.Shapes.Item(4).Top = 377.2974090748031
